# Insert a new weekly record at row 422 ("Fruta / hortaliza, semanal" update).
# This pushes the existing rows 422-496 down to 423-497 and extends the
# sheet's used range from A1:R496 to A1:R497.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("422:422").Insert()

$ws.Range("A422").Value = 11
$ws.Range("B422").Value = "Vega Monumental Concepción"
$ws.Range("C422").Value = "Bíobío"
$ws.Range("D422").Value = "2023-11-28"
$ws.Range("E422").Value = 8
$ws.Range("F422").Value = 100114013
$ws.Range("G422").Value = "Zanahoria"
$ws.Range("H422").Value = "Sin especificar"
$ws.Range("I422").Value = "Primera"
$ws.Range("J422").Value = 150
$ws.Range("K422").Value = 6500
$ws.Range("L422").Value = 6500
$ws.Range("M422").Value = 6500
$ws.Range("N422").Value = "$/saco 20 kilos"
$ws.Range("O422").Value = "Región Metropolitana"
$ws.Range("P422").Value = 325
$ws.Range("Q422").Value = 20
$ws.Range("R422").Value = "Hortaliza"
